# Updated symbol list on Fri Dec 16 08:45:55 UTC 2022 with GitHub Actions
#
# All "Price" (column D) cells in this sheet are stored as TEXT (not
# numbers), even though their content looks numeric (e.g. "262.82",
# "6.200"). To keep Excel's COM Range.Value auto-conversion from turning
# these numeric-looking strings into real numbers (which would also
# silently drop meaningful trailing zeros, e.g. "3.690" -> 3.69), every
# assignment to column D is written with a leading apostrophe, which is
# the standard Excel convention for "force text".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $ws.Range($addr).Value = "'" + $value
}

# --- simple price (column D) refreshes -----------------------------------
Set-TextValue "D2" "262.59"
Set-TextValue "D3" "24.49"
Set-TextValue "D4" "6.198"
Set-TextValue "D6" "6.746"
Set-TextValue "D7" "3.452"
Set-TextValue "D8" "1.347"
Set-TextValue "D9" "0.8014"

# --- rows 10-18: coin list shifted down by one, "One" inserted at top ----
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D10" "0.01330"
$ws.Range("E10").Value = "9OneONE"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D11" "0.1593"
$ws.Range("E11").Value = "10WazirXWRX"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D12" "0.08121"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D13" "0.03384"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D14" "0.03096"
$ws.Range("E14").Value = "13BitrueCoinBTR"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D15" "0.09342"
$ws.Range("E15").Value = "14BitMartTokenBMX"

$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D16" "3.699"
$ws.Range("E16").Value = "15MCDexMCB"

$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D17" "0.001689"
$ws.Range("E17").Value = "16BitForexTokenBF"

$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D18" "0.04811"
$ws.Range("E18").Value = "17CoinExTokenCET"

# --- remaining standalone price (column D) refreshes ----------------------
Set-TextValue "D19" "0.006209"
Set-TextValue "D20" "0.006173"
Set-TextValue "D21" "0.001102"
Set-TextValue "D22" "0.0001501"
Set-TextValue "D23" "3.690"
Set-TextValue "D26" "0.1277"
Set-TextValue "D27" "0.0006361"
Set-TextValue "D40" "0.04636"
Set-TextValue "D41" "0.007053"
Set-TextValue "D42" "0.1120"
Set-TextValue "D43" "0.003602"
Set-TextValue "D46" "0.00005882"

Set-TextValue "D49" "0.09218"
$ws.Range("E49").Value = "48BOLOBOLOWorstin24h"

Set-TextValue "D50" "0.00002101"
